# update w/ rolling + elabe and ifop poll (4/2) + ow polls of the week
#
# This script updates the PollsData sheet: corrects/refreshes the rolling,
# Elabe and Ifop poll rows (ids 198-207, previously rows 300-309) and
# appends four new poll rows (ids 208-211, rows 310-313) covering the
# OpinionWay "poll of the week", the Ifop poll from 4/2, a new Ipsos
# rolling wave, and a new Elabe poll.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 300: id 198, house opinionway, poll_type partially, regular
    $ws.Cells.Item(300, 1).Value = 198
    $ws.Cells.Item(300, 2).Value = 2022
    $ws.Cells.Item(300, 3).Value = 3
    $ws.Cells.Item(300, 4).Value = 28
    $ws.Cells.Item(300, 5).Value = 3
    $ws.Cells.Item(300, 6).Value = 29
    $ws.Cells.Item(300, 7).Value = "opinionway"
    $ws.Cells.Item(300, 8).Value = "partially"
    $ws.Cells.Item(300, 9).Value = "regular"
    $ws.Cells.Item(300, 10).Value = 772
    $ws.Cells.Item(300, 11).Value = 1
    $ws.Cells.Item(300, 12).Value = 1
    $ws.Cells.Item(300, 13).Value = 1
    $ws.Cells.Item(300, 14).Value = 1
    $ws.Cells.Item(300, 15).Value = 1
    $ws.Cells.Item(300, 16).Value = 15
    $ws.Cells.Item(300, 17).Value = 3
    $ws.Cells.Item(300, 19).Value = 5
    $ws.Cells.Item(300, 20).Value = 2
    $ws.Cells.Item(300, 22).Value = 28
    $ws.Cells.Item(300, 23).Value = 11
    $ws.Cells.Item(300, 26).Value = 3
    $ws.Cells.Item(300, 27).Value = 1
    $ws.Cells.Item(300, 28).Value = 20
    $ws.Cells.Item(300, 29).Value = 10

    # Row 301: id 199, house ipsos, poll_type excluded, rolling
    $ws.Cells.Item(301, 1).Value = 199
    $ws.Cells.Item(301, 2).Value = 2022
    $ws.Cells.Item(301, 3).Value = 3
    $ws.Cells.Item(301, 4).Value = 27
    $ws.Cells.Item(301, 5).Value = 3
    $ws.Cells.Item(301, 6).Value = 30
    $ws.Cells.Item(301, 7).Value = "ipsos"
    $ws.Cells.Item(301, 8).Value = "excluded"
    $ws.Cells.Item(301, 9).Value = "rolling"
    $ws.Cells.Item(301, 10).Value = 951
    $ws.Cells.Item(301, 11).Value = 0
    $ws.Cells.Item(301, 12).Value = 1
    $ws.Cells.Item(301, 13).Value = 0.5
    $ws.Cells.Item(301, 14).Value = 1.5
    $ws.Cells.Item(301, 15).Value = 0.5
    $ws.Cells.Item(301, 16).Value = 15.5
    $ws.Cells.Item(301, 17).Value = 3
    $ws.Cells.Item(301, 19).Value = 6
    $ws.Cells.Item(301, 20).Value = 1.5
    $ws.Cells.Item(301, 22).Value = 27
    $ws.Cells.Item(301, 23).Value = 9
    $ws.Cells.Item(301, 26).Value = 2
    $ws.Cells.Item(301, 27).Value = 2
    $ws.Cells.Item(301, 28).Value = 20.5
    $ws.Cells.Item(301, 29).Value = 11.5

    # Row 302: id 200, house elabe, poll_type partially, regular
    $ws.Cells.Item(302, 1).Value = 200
    $ws.Cells.Item(302, 2).Value = 2022
    $ws.Cells.Item(302, 3).Value = 3
    $ws.Cells.Item(302, 4).Value = 28
    $ws.Cells.Item(302, 5).Value = 3
    $ws.Cells.Item(302, 6).Value = 30
    $ws.Cells.Item(302, 7).Value = "elabe"
    $ws.Cells.Item(302, 8).Value = "partially"
    $ws.Cells.Item(302, 9).Value = "regular"
    $ws.Cells.Item(302, 10).Value = 966
    $ws.Cells.Item(302, 11).Value = 1
    $ws.Cells.Item(302, 12).Value = 1
    $ws.Cells.Item(302, 13).Value = 0.5
    $ws.Cells.Item(302, 14).Value = 1.5
    $ws.Cells.Item(302, 15).Value = 0.5
    $ws.Cells.Item(302, 16).Value = 15.5
    $ws.Cells.Item(302, 17).Value = 2.5
    $ws.Cells.Item(302, 19).Value = 4
    $ws.Cells.Item(302, 20).Value = 2
    $ws.Cells.Item(302, 22).Value = 28
    $ws.Cells.Item(302, 23).Value = 9.5
    $ws.Cells.Item(302, 26).Value = 2.5
    $ws.Cells.Item(302, 27).Value = 2.5
    $ws.Cells.Item(302, 28).Value = 21
    $ws.Cells.Item(302, 29).Value = 10.5

    # Row 303: id 201, house opinionway, poll_type partially, regular
    $ws.Cells.Item(303, 1).Value = 201
    $ws.Cells.Item(303, 2).Value = 2022
    $ws.Cells.Item(303, 3).Value = 3
    $ws.Cells.Item(303, 4).Value = 29
    $ws.Cells.Item(303, 5).Value = 3
    $ws.Cells.Item(303, 6).Value = 30
    $ws.Cells.Item(303, 7).Value = "opinionway"
    $ws.Cells.Item(303, 8).Value = "partially"
    $ws.Cells.Item(303, 9).Value = "regular"
    $ws.Cells.Item(303, 10).Value = 772
    $ws.Cells.Item(303, 11).Value = 1
    $ws.Cells.Item(303, 12).Value = 1
    $ws.Cells.Item(303, 13).Value = 1
    $ws.Cells.Item(303, 14).Value = 1
    $ws.Cells.Item(303, 15).Value = 1
    $ws.Cells.Item(303, 16).Value = 15
    $ws.Cells.Item(303, 17).Value = 3
    $ws.Cells.Item(303, 19).Value = 5
    $ws.Cells.Item(303, 20).Value = 2
    $ws.Cells.Item(303, 22).Value = 28
    $ws.Cells.Item(303, 23).Value = 11
    $ws.Cells.Item(303, 26).Value = 3
    $ws.Cells.Item(303, 27).Value = 1
    $ws.Cells.Item(303, 28).Value = 20
    $ws.Cells.Item(303, 29).Value = 10

    # Row 304: id 202, house opinionway, poll_type partially, rolling
    $ws.Cells.Item(304, 1).Value = 202
    $ws.Cells.Item(304, 2).Value = 2022
    $ws.Cells.Item(304, 3).Value = 3
    $ws.Cells.Item(304, 4).Value = 28
    $ws.Cells.Item(304, 5).Value = 3
    $ws.Cells.Item(304, 6).Value = 31
    $ws.Cells.Item(304, 7).Value = "opinionway"
    $ws.Cells.Item(304, 8).Value = "partially"
    $ws.Cells.Item(304, 9).Value = "rolling"
    $ws.Cells.Item(304, 10).Value = 1119
    $ws.Cells.Item(304, 11).Value = 1
    $ws.Cells.Item(304, 12).Value = 1
    $ws.Cells.Item(304, 13).Value = 1
    $ws.Cells.Item(304, 14).Value = 1
    $ws.Cells.Item(304, 15).Value = 1
    $ws.Cells.Item(304, 16).Value = 15
    $ws.Cells.Item(304, 17).Value = 3
    $ws.Cells.Item(304, 19).Value = 5
    $ws.Cells.Item(304, 20).Value = 2
    $ws.Cells.Item(304, 22).Value = 28
    $ws.Cells.Item(304, 23).Value = 10
    $ws.Cells.Item(304, 26).Value = 3
    $ws.Cells.Item(304, 27).Value = 2
    $ws.Cells.Item(304, 28).Value = 20
    $ws.Cells.Item(304, 29).Value = 10

    # Row 305: id 203, house ifop, poll_type included, rolling
    $ws.Cells.Item(305, 1).Value = 203
    $ws.Cells.Item(305, 2).Value = 2022
    $ws.Cells.Item(305, 3).Value = 3
    $ws.Cells.Item(305, 4).Value = 28
    $ws.Cells.Item(305, 5).Value = 3
    $ws.Cells.Item(305, 6).Value = 31
    $ws.Cells.Item(305, 7).Value = "ifop"
    $ws.Cells.Item(305, 8).Value = "included"
    $ws.Cells.Item(305, 9).Value = "rolling"
    $ws.Cells.Item(305, 10).Value = 1500
    $ws.Cells.Item(305, 11).Value = 1
    $ws.Cells.Item(305, 12).Value = 1
    $ws.Cells.Item(305, 13).Value = 0.5
    $ws.Cells.Item(305, 14).Value = 0.5
    $ws.Cells.Item(305, 15).Value = 0.5
    $ws.Cells.Item(305, 16).Value = 15.5
    $ws.Cells.Item(305, 17).Value = 4
    $ws.Cells.Item(305, 19).Value = 4.5
    $ws.Cells.Item(305, 20).Value = 1.5
    $ws.Cells.Item(305, 22).Value = 28
    $ws.Cells.Item(305, 23).Value = 10
    $ws.Cells.Item(305, 26).Value = 2
    $ws.Cells.Item(305, 27).Value = 1.5
    $ws.Cells.Item(305, 28).Value = 21
    $ws.Cells.Item(305, 29).Value = 11

    # Row 306: id 204, house cluster17, poll_type partially, regular
    $ws.Cells.Item(306, 1).Value = 204
    $ws.Cells.Item(306, 2).Value = 2022
    $ws.Cells.Item(306, 3).Value = 3
    $ws.Cells.Item(306, 4).Value = 29
    $ws.Cells.Item(306, 5).Value = 3
    $ws.Cells.Item(306, 6).Value = 31
    $ws.Cells.Item(306, 7).Value = "cluster17"
    $ws.Cells.Item(306, 8).Value = "partially"
    $ws.Cells.Item(306, 9).Value = "regular"
    $ws.Cells.Item(306, 10).Value = 2515
    $ws.Cells.Item(306, 11).Value = 0
    $ws.Cells.Item(306, 12).Value = 1
    $ws.Cells.Item(306, 13).Value = 0.5
    $ws.Cells.Item(306, 14).Value = 1
    $ws.Cells.Item(306, 15).Value = 0.5
    $ws.Cells.Item(306, 16).Value = 16
    $ws.Cells.Item(306, 17).Value = 3
    $ws.Cells.Item(306, 19).Value = 5
    $ws.Cells.Item(306, 20).Value = 2
    $ws.Cells.Item(306, 22).Value = 27
    $ws.Cells.Item(306, 23).Value = 10
    $ws.Cells.Item(306, 26).Value = 2.5
    $ws.Cells.Item(306, 27).Value = 3
    $ws.Cells.Item(306, 28).Value = 18
    $ws.Cells.Item(306, 29).Value = 12

    # Row 307: id 205, house opinionway, poll_type partially, regular
    $ws.Cells.Item(307, 1).Value = 205
    $ws.Cells.Item(307, 2).Value = 2022
    $ws.Cells.Item(307, 3).Value = 3
    $ws.Cells.Item(307, 4).Value = 30
    $ws.Cells.Item(307, 5).Value = 3
    $ws.Cells.Item(307, 6).Value = 31
    $ws.Cells.Item(307, 7).Value = "opinionway"
    $ws.Cells.Item(307, 8).Value = "partially"
    $ws.Cells.Item(307, 9).Value = "regular"
    $ws.Cells.Item(307, 10).Value = 700
    $ws.Cells.Item(307, 11).Value = 1
    $ws.Cells.Item(307, 12).Value = 1
    $ws.Cells.Item(307, 13).Value = 1
    $ws.Cells.Item(307, 14).Value = 1
    $ws.Cells.Item(307, 15).Value = 1
    $ws.Cells.Item(307, 16).Value = 15
    $ws.Cells.Item(307, 17).Value = 4
    $ws.Cells.Item(307, 19).Value = 5
    $ws.Cells.Item(307, 20).Value = 2
    $ws.Cells.Item(307, 22).Value = 27
    $ws.Cells.Item(307, 23).Value = 10
    $ws.Cells.Item(307, 26).Value = 2
    $ws.Cells.Item(307, 27).Value = 3
    $ws.Cells.Item(307, 28).Value = 21
    $ws.Cells.Item(307, 29).Value = 10

    # Row 308: id 206, house bva, poll_type partially, regular
    $ws.Cells.Item(308, 1).Value = 206
    $ws.Cells.Item(308, 2).Value = 2022
    $ws.Cells.Item(308, 3).Value = 3
    $ws.Cells.Item(308, 4).Value = 30
    $ws.Cells.Item(308, 5).Value = 3
    $ws.Cells.Item(308, 6).Value = 31
    $ws.Cells.Item(308, 7).Value = "bva"
    $ws.Cells.Item(308, 8).Value = "partially"
    $ws.Cells.Item(308, 9).Value = "regular"
    $ws.Cells.Item(308, 10).Value = 992
    $ws.Cells.Item(308, 11).Value = 0
    $ws.Cells.Item(308, 12).Value = 1
    $ws.Cells.Item(308, 13).Value = 0.5
    $ws.Cells.Item(308, 14).Value = 1
    $ws.Cells.Item(308, 15).Value = 1
    $ws.Cells.Item(308, 16).Value = 15.5
    $ws.Cells.Item(308, 17).Value = 3.5
    $ws.Cells.Item(308, 19).Value = 5
    $ws.Cells.Item(308, 20).Value = 2
    $ws.Cells.Item(308, 22).Value = 27
    $ws.Cells.Item(308, 23).Value = 9.5
    $ws.Cells.Item(308, 26).Value = 2.5
    $ws.Cells.Item(308, 27).Value = 2.5
    $ws.Cells.Item(308, 28).Value = 21
    $ws.Cells.Item(308, 29).Value = 9.5

    # Row 309: id 207, house opinionway, poll_type partially, rolling
    $ws.Cells.Item(309, 1).Value = 207
    $ws.Cells.Item(309, 2).Value = 2022
    $ws.Cells.Item(309, 3).Value = 3
    $ws.Cells.Item(309, 4).Value = 29
    $ws.Cells.Item(309, 5).Value = 4
    $ws.Cells.Item(309, 6).Value = 1
    $ws.Cells.Item(309, 7).Value = "opinionway"
    $ws.Cells.Item(309, 8).Value = "partially"
    $ws.Cells.Item(309, 9).Value = "rolling"
    $ws.Cells.Item(309, 10).Value = 1119
    $ws.Cells.Item(309, 11).Value = 1
    $ws.Cells.Item(309, 12).Value = 0.33333333333333331
    $ws.Cells.Item(309, 13).Value = 1
    $ws.Cells.Item(309, 14).Value = 1
    $ws.Cells.Item(309, 15).Value = 1
    $ws.Cells.Item(309, 16).Value = 15
    $ws.Cells.Item(309, 17).Value = 3
    $ws.Cells.Item(309, 19).Value = 5
    $ws.Cells.Item(309, 20).Value = 3
    $ws.Cells.Item(309, 22).Value = 28
    $ws.Cells.Item(309, 23).Value = 9
    $ws.Cells.Item(309, 26).Value = 3
    $ws.Cells.Item(309, 27).Value = 2
    $ws.Cells.Item(309, 28).Value = 20
    $ws.Cells.Item(309, 29).Value = 10

    # Row 310: id 208, house ifop, poll_type included, rolling
    $ws.Cells.Item(310, 1).Value = 208
    $ws.Cells.Item(310, 2).Value = 2022
    $ws.Cells.Item(310, 3).Value = 3
    $ws.Cells.Item(310, 4).Value = 29
    $ws.Cells.Item(310, 5).Value = 4
    $ws.Cells.Item(310, 6).Value = 1
    $ws.Cells.Item(310, 7).Value = "ifop"
    $ws.Cells.Item(310, 8).Value = "included"
    $ws.Cells.Item(310, 9).Value = "rolling"
    $ws.Cells.Item(310, 10).Value = 1500
    $ws.Cells.Item(310, 11).Value = 1
    $ws.Cells.Item(310, 12).Value = 0.33333333333333331
    $ws.Cells.Item(310, 13).Value = 0.5
    $ws.Cells.Item(310, 14).Value = 1
    $ws.Cells.Item(310, 15).Value = 0.5
    $ws.Cells.Item(310, 16).Value = 15
    $ws.Cells.Item(310, 17).Value = 3.5
    $ws.Cells.Item(310, 19).Value = 4.5
    $ws.Cells.Item(310, 20).Value = 1.5
    $ws.Cells.Item(310, 22).Value = 28
    $ws.Cells.Item(310, 23).Value = 9.5
    $ws.Cells.Item(310, 26).Value = 2
    $ws.Cells.Item(310, 27).Value = 2
    $ws.Cells.Item(310, 28).Value = 21.5
    $ws.Cells.Item(310, 29).Value = 11

    # Row 311: id 209, house ifop, poll_type included, regular
    $ws.Cells.Item(311, 1).Value = 209
    $ws.Cells.Item(311, 2).Value = 2022
    $ws.Cells.Item(311, 3).Value = 3
    $ws.Cells.Item(311, 4).Value = 31
    $ws.Cells.Item(311, 5).Value = 4
    $ws.Cells.Item(311, 6).Value = 1
    $ws.Cells.Item(311, 7).Value = "ifop"
    $ws.Cells.Item(311, 8).Value = "included"
    $ws.Cells.Item(311, 9).Value = "regular"
    $ws.Cells.Item(311, 10).Value = 1000
    $ws.Cells.Item(311, 11).Value = 1
    $ws.Cells.Item(311, 12).Value = 1
    $ws.Cells.Item(311, 13).Value = 0.5
    $ws.Cells.Item(311, 14).Value = 1
    $ws.Cells.Item(311, 15).Value = 0.5
    $ws.Cells.Item(311, 16).Value = 15
    $ws.Cells.Item(311, 17).Value = 3.5
    $ws.Cells.Item(311, 19).Value = 5
    $ws.Cells.Item(311, 20).Value = 2
    $ws.Cells.Item(311, 22).Value = 27
    $ws.Cells.Item(311, 23).Value = 9
    $ws.Cells.Item(311, 26).Value = 2.5
    $ws.Cells.Item(311, 27).Value = 2
    $ws.Cells.Item(311, 28).Value = 22
    $ws.Cells.Item(311, 29).Value = 10.5

    # Row 312: id 210, house ipsos, poll_type excluded, rolling
    $ws.Cells.Item(312, 1).Value = 210
    $ws.Cells.Item(312, 2).Value = 2022
    $ws.Cells.Item(312, 3).Value = 3
    $ws.Cells.Item(312, 4).Value = 30
    $ws.Cells.Item(312, 5).Value = 4
    $ws.Cells.Item(312, 6).Value = 2
    $ws.Cells.Item(312, 7).Value = "ipsos"
    $ws.Cells.Item(312, 8).Value = "excluded"
    $ws.Cells.Item(312, 9).Value = "rolling"
    $ws.Cells.Item(312, 10).Value = 1066
    $ws.Cells.Item(312, 11).Value = 0
    $ws.Cells.Item(312, 12).Value = 1
    $ws.Cells.Item(312, 13).Value = 0.5
    $ws.Cells.Item(312, 14).Value = 1.5
    $ws.Cells.Item(312, 15).Value = 0.5
    $ws.Cells.Item(312, 16).Value = 15.5
    $ws.Cells.Item(312, 17).Value = 3
    $ws.Cells.Item(312, 19).Value = 6
    $ws.Cells.Item(312, 20).Value = 2
    $ws.Cells.Item(312, 22).Value = 26
    $ws.Cells.Item(312, 23).Value = 9.5
    $ws.Cells.Item(312, 26).Value = 2.5
    $ws.Cells.Item(312, 27).Value = 1.5
    $ws.Cells.Item(312, 28).Value = 21
    $ws.Cells.Item(312, 29).Value = 11

    # Row 313: id 211, house elabe, poll_type partially, regular
    $ws.Cells.Item(313, 1).Value = 211
    $ws.Cells.Item(313, 2).Value = 2022
    $ws.Cells.Item(313, 3).Value = 3
    $ws.Cells.Item(313, 4).Value = 31
    $ws.Cells.Item(313, 5).Value = 4
    $ws.Cells.Item(313, 6).Value = 2
    $ws.Cells.Item(313, 7).Value = "elabe"
    $ws.Cells.Item(313, 8).Value = "partially"
    $ws.Cells.Item(313, 9).Value = "regular"
    $ws.Cells.Item(313, 10).Value = 999
    $ws.Cells.Item(313, 11).Value = 0
    $ws.Cells.Item(313, 12).Value = 1
    $ws.Cells.Item(313, 13).Value = 0.5
    $ws.Cells.Item(313, 14).Value = 1.5
    $ws.Cells.Item(313, 15).Value = 0.5
    $ws.Cells.Item(313, 16).Value = 15
    $ws.Cells.Item(313, 17).Value = 3.5
    $ws.Cells.Item(313, 19).Value = 4.5
    $ws.Cells.Item(313, 20).Value = 1.5
    $ws.Cells.Item(313, 22).Value = 28.5
    $ws.Cells.Item(313, 23).Value = 8.5
    $ws.Cells.Item(313, 26).Value = 2.5
    $ws.Cells.Item(313, 27).Value = 2.5
    $ws.Cells.Item(313, 28).Value = 22
    $ws.Cells.Item(313, 29).Value = 9.5

    # Two additional blank placeholder rows appended after the data block,
    # matching the style used throughout column L/M (rounding helper column)
    $ws.Cells.Item(463, 12).NumberFormat = "0.00"
    $ws.Cells.Item(463, 13).NumberFormat = "0.00"
    $ws.Cells.Item(464, 12).NumberFormat = "0.00"
    $ws.Cells.Item(464, 13).NumberFormat = "0.00"

    # Restore selection to the last touched cell of the data block
    $ws.Range("AD311").Select()
